# Add a "wall" row/column to the collision table and fill in its values,
# per the commit: "Bug fixes + knockback + Start on walls".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared string / header label used for both the new row and new column.
$wallLabel = "wall"

# --- New column L (wall) -------------------------------------------------
$ws.Range("L1").Value = $wallLabel

# Column L values for existing rows 2-11 (wall vs. each actor):
# no-collision ("▬") for player/foes/friendbullet, collision ("x") for the
# enemy bullets, lostsoul and altar - mirrors the matching new row below.
$ws.Range("L2").Value  = "▬"
$ws.Range("L3").Value  = "▬"
$ws.Range("L4").Value  = "▬"
$ws.Range("L5").Value  = "▬"
$ws.Range("L6").Value  = "▬"
$ws.Range("L7").Value  = "▬"
$ws.Range("L8").Value  = "▬"
$ws.Range("L9").Value  = "▬"
$ws.Range("L10").Value = "x"
$ws.Range("L11").Value = "x"

# --- New row 12 (wall) ----------------------------------------------------
$ws.Range("A12").Value = $wallLabel
$ws.Range("B12").Value = "▬"
$ws.Range("C12").Value = "▬"
$ws.Range("D12").Value = "▬"
$ws.Range("E12").Value = "▬"
$ws.Range("F12").Value = "▬"
$ws.Range("G12").Value = "▬"
$ws.Range("H12").Value = "▬"
$ws.Range("I12").Value = "▬"
$ws.Range("J12").Value = "x"
$ws.Range("K12").Value = "x"
$ws.Range("L12").Value = "x"

# Match the author's final selection/cursor position recorded in the diff.
$ws.Range("L11").Select()
